# Generate Report for Handback
#
# A new handback attempt was recorded for 6e1ba82c-4df3-4bc3-9576-2ae31e46d890
# in both the zh-cn and de-de target-language sheets: the handback file that
# came in was not the latest revision, so a second "Latest Target
# File"/"Latest Handback File"/"Latest Handback DateTime" entry is recorded
# in columns I/J/K of row 6, and the mismatch is reported in column P
# ("Error Detail"). The Error Detail column is also widened so the message
# is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cd56859e8f6a5bea7ddb65d12b888e86fbf61e3/e2e/6e1ba82c-4df3-4bc3-9576-2ae31e46d890.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/332471f68a8349344cac028f01e3d15d7198b9c9/e2e/6e1ba82c-4df3-4bc3-9576-2ae31e46d890.md."

$sheets = @{
    "zh-cn" = @{
        HandoffFile = "6e1ba82c-4df3-4bc3-9576-2ae31e46d890.fa97473eacb4fcc367cdcd29df3135cfa4247a87.zh-cn.xlf"
        HandbackDateTime = "2016-09-03 18:55:57"
    }
    "de-de" = @{
        HandoffFile = "6e1ba82c-4df3-4bc3-9576-2ae31e46d890.fa97473eacb4fcc367cdcd29df3135cfa4247a87.de-de.xlf"
        HandbackDateTime = "2016-09-03 18:56:09"
    }
}

foreach ($sheetName in $sheets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheets[$sheetName]

    # Latest Target File (I6) - mirrors A6's handback-to-md hyperlink.
    $ws.Range("I6").Value = "6e1ba82c-4df3-4bc3-9576-2ae31e46d890.md"
    $ws.Range("I6").Style = "Hyperlink"
    $ws.Hyperlinks.Add(
        $ws.Range("I6"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/332471f68a8349344cac028f01e3d15d7198b9c9/e2e/6e1ba82c-4df3-4bc3-9576-2ae31e46d890.md",
        "",
        "",
        "6e1ba82c-4df3-4bc3-9576-2ae31e46d890.md"
    ) | Out-Null

    # Latest Handback File (J6) - same xlf as the Latest Handoff File (G6).
    $ws.Range("J6").Value = $info.HandoffFile

    # Latest Handback DateTime (K6).
    $ws.Range("K6").Value = $info.HandbackDateTime

    # Error Detail (P6) - explains why this handback triggered a new row.
    $ws.Range("P6").Value = $errorDetail

    # Error Detail column (P) is widened to fit the long message.
    $ws.Range("P1").ColumnWidth = 39.15
}
